$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# --- Grow the table from A1:C2 to A1:D15, adding a 4th column ---
$lo.Resize($ws.Range("A1:D2"))
$ws.Range("D1").Value = "Lead"

# Move the new "Lead" column so it sits before "Message" (C <-> D swap),
# i.e. final order: Date, Bible Verse, Lead, Message
$ws.Columns.Item(4).Cut()
$ws.Columns.Item(3).Insert()

# Now resize the table down to the full 15 rows
$lo.Resize($ws.Range("A1:D15"))

# --- Fill in the weekly dates in column A ---
$ws.Range("A2").Value = 45925
$ws.Range("A3").Formula = "=A2+7"
$ws.Range("A4:A15").Formula = "=A3+7"

$wb.Worksheets.Item(1).Activate()
$ws.Range("C2").Select()
